$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last month block (rows 48:59 - divider, blank spacer,
# header, 6 activity rows, total row, Paid row, Not Paid row) into the
# new block at rows 60:71 so formatting/styles match exactly.
$ws.Range("A48:G59").Copy($ws.Range("A60"))

# -- Row 62: new month header ("Aban 98") --
$ws.Range("A62").Value = "آبان 98"

# -- Rows 63-68: activity rows for the new month --
# Row 63: * BronchoVision GUI -> no hours yet, task note updated
$ws.Range("C63").ClearContents()
$ws.Range("E63").Value = "• Fill holes in 3D view"

# Row 64: * 2D Views -> no hours yet, clear old task note
$ws.Range("C64").ClearContents()
$ws.Range("E64").ClearContents()

# Row 65: * 3D View & Virtual Camera -> 2 hours, clear old task note
$ws.Range("C65").Value = 2
$ws.Range("E65").ClearContents()

# Row 66: * Tracker connection -> no hours yet, clear old task note
$ws.Range("C66").ClearContents()
$ws.Range("E66").ClearContents()

# Row 67: * Data Importing -> no hours yet (E67 already blank)
$ws.Range("C67").ClearContents()

# Row 68: * Multi-threading -> no hours yet (E68 already blank)
$ws.Range("C68").ClearContents()

# -- Row 69: Total Hours formula, re-anchored to the new block --
$ws.Range("C69").Formula = "=SUM(C63:C68)"

# -- Rows 70-71: Paid / Not Paid hours for the new month --
$ws.Range("D70").ClearContents()
$ws.Range("D71").Value = 2

# Update the active selection to match the edited cell.
$ws.Range("C63").Select()
